{"js": "// Office.js (Word JavaScript API) port of the commit that re-rolls the 100\n// addition/subtraction \"within 100\" drill problems laid out in the\n// document's single 20-row x 5-column table. The document's first\n// paragraph (the date heading) is untouched; only the table cell text is\n// replaced, in row-major reading order, cell formatting (TimeNewRoman,\n// sz 30) is left exactly as-is.\nconst newValues = [\n  \"72-6=\", \"26+41=\", \"50+29=\", \"27+11=\", \"75+3=\", \"83-36=\", \"84-55=\", \"73-10=\", \"3+72=\", \"66+8=\",\n  \"44-5=\", \"10+44=\", \"4+54=\", \"41+8=\", \"4+42=\", \"19-12=\", \"22+45=\", \"77+14=\", \"21-7=\", \"72-37=\",\n  \"37+4=\", \"46+50=\", \"17+16=\", \"33-0=\", \"78-70=\", \"71-28=\", \"23+1=\", \"71+3=\", \"68-55=\", \"28+64=\",\n  \"31+52=\", \"91-38=\", \"4+82=\", \"25-12=\", \"3+42=\", \"3+68=\", \"42-36=\", \"39-19=\", \"84-34=\", \"55-43=\",\n  \"55+16=\", \"36-28=\", \"84-82=\", \"59+3=\", \"37+58=\", \"54-10=\", \"22-19=\", \"59-59=\", \"81-35=\", \"64-25=\",\n  \"25+44=\", \"50-33=\", \"19+67=\", \"39+3=\", \"61-14=\", \"74-20=\", \"91-8=\", \"56+10=\", \"23-1=\", \"15+38=\",\n  \"89-40=\", \"41+36=\", \"24+36=\", \"90-22=\", \"50-47=\", \"4+7=\", \"13+27=\", \"47+4=\", \"43-5=\", \"97-73=\",\n  \"64+18=\", \"57-12=\", \"27-24=\", \"66+9=\", \"77-57=\", \"80-62=\", \"26+43=\", \"7+17=\", \"96-92=\", \"79-23=\",\n  \"78-39=\", \"48+49=\", \"31+42=\", \"26+20=\", \"1+82=\", \"90-34=\", \"0+31=\", \"90-18=\", \"16+79=\", \"16+38=\",\n  \"55-54=\", \"20+72=\", \"82+4=\", \"87-18=\", \"89+1=\", \"29+65=\", \"50+44=\", \"97-17=\", \"90-72=\", \"94-22=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\nconst table = tables.items[0];\n// `columnCount` isn't a loadable property in this host, so derive the grid\n// shape from `values` (a 2-D array mirroring the table's rows/columns).\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst columnCount = table.values.length > 0 ? table.values[0].length : 0;\n\nif (rowCount * columnCount !== newValues.length) {\n  throw new Error(\n    `Table shape ${rowCount}x${columnCount} does not hold ${newValues.length} values.`\n  );\n}\n\n// Grab the (single) paragraph inside each cell, in row-major order, which\n// is the same order the new values/commit diff enumerate the cells in.\nconst paragraphs = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    paragraphs.push(cell.body.paragraphs.getFirst());\n  }\n}\nawait context.sync();\n\n// Replace the text through each paragraph's Range rather than the cell\n// body, so the existing run (font/size) and paragraph properties survive\n// untouched \u2014 only the literal characters change, matching the diff.\nfor (let i = 0; i < paragraphs.length; i++) {\n  paragraphs[i].getRange().insertText(newValues[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word COM-interop (PowerShell) port of the commit that re-rolls the 100\n# addition/subtraction \"within 100\" drill problems laid out in the\n# document's single 20-row x 5-column table. The document's first\n# paragraph (the date heading) is untouched; only the table cell text is\n# replaced, in row-major reading order, and cell formatting\n# (TimeNewRoman, sz 30) is left exactly as-is.\n$newValues = @(\n    \"72-6=\", \"26+41=\", \"50+29=\", \"27+11=\", \"75+3=\", \"83-36=\", \"84-55=\", \"73-10=\", \"3+72=\", \"66+8=\",\n    \"44-5=\", \"10+44=\", \"4+54=\", \"41+8=\", \"4+42=\", \"19-12=\", \"22+45=\", \"77+14=\", \"21-7=\", \"72-37=\",\n    \"37+4=\", \"46+50=\", \"17+16=\", \"33-0=\", \"78-70=\", \"71-28=\", \"23+1=\", \"71+3=\", \"68-55=\", \"28+64=\",\n    \"31+52=\", \"91-38=\", \"4+82=\", \"25-12=\", \"3+42=\", \"3+68=\", \"42-36=\", \"39-19=\", \"84-34=\", \"55-43=\",\n    \"55+16=\", \"36-28=\", \"84-82=\", \"59+3=\", \"37+58=\", \"54-10=\", \"22-19=\", \"59-59=\", \"81-35=\", \"64-25=\",\n    \"25+44=\", \"50-33=\", \"19+67=\", \"39+3=\", \"61-14=\", \"74-20=\", \"91-8=\", \"56+10=\", \"23-1=\", \"15+38=\",\n    \"89-40=\", \"41+36=\", \"24+36=\", \"90-22=\", \"50-47=\", \"4+7=\", \"13+27=\", \"47+4=\", \"43-5=\", \"97-73=\",\n    \"64+18=\", \"57-12=\", \"27-24=\", \"66+9=\", \"77-57=\", \"80-62=\", \"26+43=\", \"7+17=\", \"96-92=\", \"79-23=\",\n    \"78-39=\", \"48+49=\", \"31+42=\", \"26+20=\", \"1+82=\", \"90-34=\", \"0+31=\", \"90-18=\", \"16+79=\", \"16+38=\",\n    \"55-54=\", \"20+72=\", \"82+4=\", \"87-18=\", \"89+1=\", \"29+65=\", \"50+44=\", \"97-17=\", \"90-72=\", \"94-22=\"\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables(1)\n$rows = $table.Rows.Count\n$cols = $table.Columns.Count\n\nif (($rows * $cols) -ne $newValues.Count) {\n    throw \"Table shape $rows x $cols does not hold $($newValues.Count) values.\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cellRange = $cell.Range\n        # A cell Range includes the trailing end-of-cell marker character;\n        # back the end up one character so only the visible formula text\n        # is overwritten. Assigning .Text keeps the run's existing\n        # character formatting (font, size) because Word reuses the first\n        # run's properties for the replacement text.\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        $cellRange.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
